# "Changed cleaning to support additional data"
# The worksheet holds a 2-column lookup list (A = numeric code, B = id string).
# The cleaning/sort order of the existing 7 rows changes, and 38 brand-new
# rows of data are appended, growing the sheet from A1:B8 to A1:B46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full desired contents for rows 2..46 (column A value, column B id), in order.
$data = @(
    @(12, "N931325309008"),
    @(51, "W931252509017"),
    @(34, "T931252911047"),
    @(38, "T931252110004"),
    @(18, "N931253409013"),
    @(15, "V931414517045"),
    @(31, "R931100609011"),
    @(28, "R931321009045"),
    @(32, "R931412017031"),
    @(37, "L931412020028"),
    @(39, "L931321113001"),
    @(17, "D931383810007"),
    @(30, "V931252909047"),
    @(46, "Z931100609006"),
    @(45, "H931101008036"),
    @(47, "Y931412017035"),
    @(48, "C931253110015"),
    @(36, "V931240110042"),
    @(35, "H931321309010"),
    @(49, "R931100609009"),
    @(33, "X931252710015"),
    @(50, "N931240110007"),
    @(42, "T931100609002"),
    @(40, "M931321110016"),
    @(0, "G936239910030"),
    @(52, "L931101008038"),
    @(2, "M888201710014"),
    @(3, "D931252109051"),
    @(5, "V802229210007"),
    @(6, "D931100608056"),
    @(7, "C931316110004"),
    @(8, "A931259308039"),
    @(13, "G931259509014"),
    @(14, "J931101108070"),
    @(16, "K931100609063"),
    @(19, "G931321110032"),
    @(20, "W931321110033"),
    @(21, "V931412017033"),
    @(22, "V931321008075"),
    @(23, "U931325209009"),
    @(24, "W310350110023"),
    @(25, "W931254310067"),
    @(26, "F931100509027"),
    @(29, "D931325309031"),
    @(53, "B931400418001")
)

$startRow = 2
$lastRow = $startRow + $data.Length - 1

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Row 2 already carries the correct formatting (bold/bordered/centered style
# for column A, plain for column B); stamp that same formatting onto every
# newly added row so the new rows 9-46 look like rows 2-8 did before.
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

